# Update with latest cht-conf changes again and remove NO_LABEL
#
# The "survey" sheet has two cells (C3 and C10) that hold a literal
# "NO_LABEL" placeholder in the label column. The upstream cht-conf
# tooling no longer emits that placeholder, so both cells are removed
# entirely (not just blanked) to match the regenerated form.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

$ws.Range("C3").Clear()
$ws.Range("C10").Clear()

# The saved cursor position on the "survey" sheet also moved (it now
# rests on the second now-unused NO_LABEL row, C10, instead of C17).
$ws.Range("C10").Select() | Out-Null
